$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values per row, reflecting repulled/pushed data and mean calculation
$values = @{
    2  = 2
    3  = 1
    4  = -5
    5  = 4
    7  = 2
    8  = 2
    9  = -1
    10 = 1
    11 = -2
    12 = 9
    13 = 2
    14 = 2
    15 = 1
    16 = 1
    17 = -3
    18 = 2
    19 = 4
    20 = 3
    21 = 2
    22 = 2
    23 = 5
    24 = 1
    25 = 1
    26 = -4
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 6).Value = $values[$row]
}
